# Generate Report for Handoff
#
# The "b.md" row (row 3) on the Overview / zh-cn / de-de sheets is refreshed
# to reflect a new handoff pass: status flips from "Handed back: in sync
# with en-US" to "Ready for handoff", a new handoff xliff file + timestamp
# is recorded, "Content Duplicate" drops back to False, and an error detail
# message is attached explaining the stale handback file. The "Error Detail"
# column is also widened so the longer message is readable.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e38fcdf23a583ac96879b409072ab66abcbcc5b8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87cb20578a3d6f46d9a9eb1298646de4f9297959/e2e/b.md."

# Excel's ColumnWidth (character units) is offset from the stored OOXML
# <col width> by the default column padding (5/6 of a character here), so
# asking for an on-disk width of 40 means setting ColumnWidth to 40 - 5/6.
$targetColWidth40 = 40 - (5 / 6)

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = "2016-10-27 08:07:46"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
# Reuse the existing text "False" cell (F2) so the written value stays a
# text/shared-string cell instead of Excel auto-coercing the literal
# string "False" into a native boolean.
$wsZhCn.Range("F2").Copy($wsZhCn.Range("F3"))
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-10-27 08:07:33"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColWidth40

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("F2").Copy($wsDeDe.Range("F3"))
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-10-27 08:07:46"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColWidth40
